# Apply new predicted prices and compared returns to the S&P500 returns
# compared annual sheet. Columns G (Return_with_prediction) and H
# (return_pct_change) are updated for rows 2-29, and I2
# (mean_return_pct_change) is updated as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = @{ G = 0.03824976438307783;  H = -20.6940248254145 }
    3  = @{ G = 0.04483339027015782;  H = 16.88571514134753 }
    4  = @{ G = -0.4585870876927937;  H = -1.285227615563234 }
    5  = @{ G = -0.4388679390381338;  H = 8.346725151898845 }
    6  = @{ G = 0.2396970430948253;   H = 2.59613466691876 }
    7  = @{ G = 0.2448388215796232;   H = 11.00045457604885 }
    8  = @{ G = 0.1760996157766723;   H = 5.571234362668513 }
    9  = @{ G = 0.1703927428324036;   H = -0.9409769255812588 }
    10 = @{ G = -0.003724392634685596; H = 21.41772547748926 }
    11 = @{ G = -0.004627496247610453; H = 68.42440882121728 }
    12 = @{ G = 0.1400234185610019;   H = 2.414488481887021 }
    13 = @{ G = 0.1401382182428771;   H = 12.4365667464733 }
    14 = @{ G = 0.2667501122088802;   H = 7.850387641383325 }
    15 = @{ G = 0.2650433889716831;   H = 4.894144909519592 }
    16 = @{ G = 0.1452364219830415;   H = -5.364493435279531 }
    17 = @{ G = 0.144694201943753;    H = -4.185669076198677 }
    18 = @{ G = -0.01136787339824208; H = 30.57747494092096 }
    19 = @{ G = -0.009750159221883478; H = -1057.888263618506 }
    20 = @{ G = 0.1494298754285556;   H = 7.765421753616149 }
    21 = @{ G = 0.138699445637133;    H = -3.071418060385385 }
    22 = @{ G = 0.1674259144532954;   H = -10.09114773460152 }
    23 = @{ G = 0.1789295658707026;   H = -0.2958524821361691 }
    24 = @{ G = -0.1018184128685597; H = -7.861615868659301 }
    25 = @{ G = -0.1070734714531105; H = -7.494360711349659 }
    26 = @{ G = 0.2383632479720952;   H = 3.572951045547969 }
    27 = @{ G = 0.2283790264084297;   H = -1.802971781993691 }
    28 = @{ G = 0.0712192754139809;   H = 21.11449632731264 }
    29 = @{ G = 0.08205040369027952;  H = 16.24418789111791 }
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row].G
    $ws.Range("H$row").Value = $newValues[$row].H
}

$ws.Range("I2").Value = -31.42369657858217
